# Refresh crypto price/volume data (GitHub Actions scheduled update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '59.040.62'
$ws.Range("E2").Value = '  +2.39%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '2.589.79'
$ws.Range("E3").Value = '  +1.12%  '

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.27'
$ws.Range("E5").Value = '  +3.76%  '

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.60'
$ws.Range("E6").Value = '  +0.16%  '

# Row 7: update E7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  +1.08%  '

# Row 9: update D9, E9
$ws.Range("D9").Value = '2.602.55'
$ws.Range("E9").Value = '  +1.19%  '

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.45'
$ws.Range("E10").Value = '  +0.48%  '

# Row 11: update E11
$ws.Range("E11").Value = '  +2.25%  '

# Row 12: update E12
$ws.Range("E12").Value = '  +0.84%  '

# Row 13: update E13
$ws.Range("E13").Value = '  +3.82%  '

# Row 14: update D14, E14
$ws.Range("D14").Value = '3.045.47'
$ws.Range("E14").Value = '  +0.98%  '

# Row 15: update D15, E15
$ws.Range("D15").Value = '58.983.25'
$ws.Range("E15").Value = '  +2.28%  '

# Row 16: update E16
$ws.Range("E16").Value = '  +2.57%  '

# Row 17: update D17, E17
$ws.Range("D17").Value = '2.629.77'
$ws.Range("E17").Value = '  +2.21%  '

# Row 18: update E18
$ws.Range("E18").Value = '  +1.42%  '

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.92'
$ws.Range("E19").Value = '  +4.22%  '

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.33'
$ws.Range("E20").Value = '  +2.22%  '

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.09'
$ws.Range("E21").Value = '  +0.95%  '

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.44'
$ws.Range("E22").Value = '  +2.26%  '

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.10%  '

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.37'
$ws.Range("E24").Value = '  +4.03%  '

# Row 25: update E25
$ws.Range("E25").Value = '  +0.59%  '

# Row 26: update D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.405'
$ws.Range("E26").Value = '  +2.61%  '

# Row 27: update E27
$ws.Range("E27").Value = '  -0.03%  '

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.11'
$ws.Range("E28").Value = '  +2.98%  '

# Row 29: update E29
$ws.Range("E29").Value = '  -0.07%  '

# Row 30: update D30, E30
$ws.Range("D30").Value = '0.0₃0723'
$ws.Range("E30").Value = '  +1.37%  '

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.61'
$ws.Range("E31").Value = '  +3.95%  '

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.89'
$ws.Range("E32").Value = '  -1.82%  '

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.72'
$ws.Range("E33").Value = '  +1.57%  '

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.92'
$ws.Range("E34").Value = '  +0.01%  '

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("E35").Value = '  +1.95%  '

# Row 36: update E36
$ws.Range("E36").Value = '  +1.05%  '

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.81'
$ws.Range("E37").Value = '  +3.36%  '

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.48'
$ws.Range("E38").Value = '  +3.71%  '

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.825'
$ws.Range("E39").Value = '  -0.31%  '

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("E40").Value = '  -0.37%  '

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.53'
$ws.Range("E41").Value = '  +2.54%  '

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.18%  '

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.78'
$ws.Range("E43").Value = '  +0.97%  '

# Row 44: update B44, C44, D44, E44
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '269.27'
$ws.Range("E44").Value = '  +1.35%  '

# Row 45: update B45, C45, D45, E45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.596'
$ws.Range("E45").Value = '  +0.67%  '

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0956'
$ws.Range("E46").Value = '  +1.78%  '

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0517'
$ws.Range("E47").Value = '  +0.87%  '

# Row 48: update D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.36'
$ws.Range("E48").Value = '  +0.18%  '

# Row 49: update D49, E49
$ws.Range("D49").Value = '1.957.69'
$ws.Range("E49").Value = '  +0.75%  '

# Row 50: update E50
$ws.Range("E50").Value = '  +2.15%  '

# Row 51: update B51, C51, D51, E51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.53'
$ws.Range("E51").Value = '  +2.39%  '

